$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '45.800.15'
$ws.Range("E2").Value = '  +7.35%  '
$ws.Range("D3").Value = '2.399.63'
$ws.Range("E3").Value = '  +4.28%  '
$ws.Range("E4").Value = '  +0.20%  '
$ws.Range("D5").Value = '''113.60'
$ws.Range("E5").Value = '  +8.36%  '
$ws.Range("D6").Value = '''318.52'
$ws.Range("E6").Value = '  +2.54%  '
$ws.Range("D7").Value = '''0.635'
$ws.Range("E7").Value = '  +1.38%  '
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("E9").Value = '  +3.74%  '
$ws.Range("D10").Value = '''41.91'
$ws.Range("E10").Value = '  +5.66%  '
$ws.Range("D11").Value = '''0.0932'
$ws.Range("E11").Value = '  +3.08%  '
$ws.Range("E12").Value = '  +5.28%  '
$ws.Range("E13").Value = '  +2.26%  '
$ws.Range("E14").Value = '  +1.50%  '
$ws.Range("D15").Value = '''15.81'
$ws.Range("E15").Value = '  +3.57%  '
$ws.Range("D16").Value = '2.764.72'
$ws.Range("E16").Value = '  +4.29%  '
$ws.Range("D17").Value = '2.394.20'
$ws.Range("E17").Value = '  +4.11%  '
$ws.Range("D18").Value = '45.764.13'
$ws.Range("E18").Value = '  +7.31%  '
$ws.Range("D19").Value = '''7.48'
$ws.Range("E19").Value = '  +2.22%  '
$ws.Range("D20").Value = '''0.0000109'
$ws.Range("E20").Value = '  +3.45%  '
$ws.Range("D21").Value = '''13.39'
$ws.Range("E21").Value = '  -2.05%  '
$ws.Range("D22").Value = '''74.57'
$ws.Range("E22").Value = '  +1.38%  '
$ws.Range("E23").Value = '  +1.88%  '
$ws.Range("D24").Value = '''264.58'
$ws.Range("E24").Value = '  -0.82%  '
$ws.Range("D25").Value = '''2.34'
$ws.Range("E25").Value = '  +4.48%  '
$ws.Range("E26").Value = '  -0.59%  '
$ws.Range("B27").Value = 'Filecoin'
$ws.Range("C27").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D27").Value = '''7.64'
$ws.Range("E27").Value = '  +1.18%  '
$ws.Range("B28").Value = 'Cosmos'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D28").Value = '''11.31'
$ws.Range("E28").Value = '  +3.19%  '
$ws.Range("E29").Value = '  +2.53%  '
$ws.Range("D30").Value = '''39.01'
$ws.Range("E30").Value = '  +2.25%  '
$ws.Range("D31").Value = '''22.79'
$ws.Range("E31").Value = '  +2.45%  '
$ws.Range("D32").Value = '''0.0978'
$ws.Range("E32").Value = '  +12.84%  '
$ws.Range("D33").Value = '''172.35'
$ws.Range("E33").Value = '  +4.43%  '
$ws.Range("D34").Value = '''2.96'
$ws.Range("E34").Value = '  +4.49%  '
$ws.Range("D35").Value = '''0.133'
$ws.Range("E35").Value = '  +1.72%  '
$ws.Range("D36").Value = '''0.119'
$ws.Range("E36").Value = '  +5.69%  '
$ws.Range("D37").Value = '''4.90'
$ws.Range("E37").Value = '  +5.81%  '
$ws.Range("D38").Value = '''4.10'
$ws.Range("E38").Value = '  +13.30%  '
$ws.Range("D39").Value = '''3.04'
$ws.Range("E39").Value = '  +8.13%  '
$ws.Range("E40").Value = '  +1.92%  '
$ws.Range("D41").Value = '''1.79'
$ws.Range("E41").Value = '  +14.22%  '
$ws.Range("D42").Value = '''102.27'
$ws.Range("E42").Value = '  -5.04%  '
$ws.Range("D43").Value = '''0.240'
$ws.Range("E43").Value = '  +5.08%  '
$ws.Range("E44").Value = '  +9.93%  '
$ws.Range("D45").Value = '''71.69'
$ws.Range("E45").Value = '  +0.10%  '
$ws.Range("D46").Value = '''87.97'
$ws.Range("E46").Value = '  +15.46%  '
$ws.Range("E47").Value = '  +0.18%  '
$ws.Range("D48").Value = '''115.48'
$ws.Range("E48").Value = '  +3.45%  '
$ws.Range("B49").Value = 'FraxShare'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D49").Value = '''9.50'
$ws.Range("E49").Value = '  +6.99%  '
$ws.Range("B50").Value = 'THORChain'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D50").Value = '''5.62'
$ws.Range("E50").Value = '  +8.58%  '
$ws.Range("D51").Value = '1.660.55'
$ws.Range("E51").Value = '  -2.70%  '
